$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (Transect number / size columns)
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0

# Update the active cell / selection on the frozen (bottom-left) pane to A2
$ws.Range("A2").Select()
